$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1 - match style of existing header cells (e.g. G1)
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Save column values for rows 2-11 (1 = save, 0 = no save)
$saveValues = @(0, 0, 0, 0, 1, 0, 0, 0, 0, 0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Range("H$row").Value = $saveValues[$i]
}
